$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix A126:A132 from text "21" to numeric 21 ---
for ($r = 126; $r -le 132; $r++) {
    $ws.Cells.Item($r, 1).Value2 = 21
}

# --- New row 133 data ---
$D133 = @'
Looking for someone that can do AI dubbed over voice video production from existing reels 
Budget
: $100
Posted On
: June 15, 2024 09:04 UTC
Category
: Video Editing
Skills
:Adobe After Effects,     Video Editing,     Video Post-Editing    
Skills
:        Adobe After Effects,                     Video Editing,                     Video Post-Editing            
Country
: Spain
click to apply

'@
$E133 = @'
Looking for someone that can do AI dubbed over voice video production from existing reels <br /><br /><b>Budget</b>: $100
<br /><b>Posted On</b>: June 15, 2024 09:04 UTC<br /><b>Category</b>: Video Editing<br /><b>Skills</b>:Adobe After Effects,     Video Editing,     Video Post-Editing    
<br /><b>Skills</b>:        Adobe After Effects,                     Video Editing,                     Video Post-Editing            <br /><b>Country</b>: Spain
<br /><a href="https://www.upwork.com/jobs/mudded-video-editing_%7E01f58b6b0d38bf9a23?source=rss">click to apply</a>

'@

# Column A: text "23" (number-looking text -> force text format)
$ws.Cells.Item(133, 1).NumberFormat = "@"
$ws.Cells.Item(133, 1).Value2 = "23"
$ws.Cells.Item(133, 2).Value2 = "Ai mudded video editing  - Upwork"
$ws.Cells.Item(133, 3).Value2 = "https://www.upwork.com/jobs/mudded-video-editing_%7E01f58b6b0d38bf9a23?source=rss"
$ws.Cells.Item(133, 4).Value2 = $D133
$ws.Cells.Item(133, 5).Value2 = $E133
$ws.Cells.Item(133, 6).Value2 = "Sat, 15 Jun 2024 09:04:36 +0000"
$ws.Cells.Item(133, 7).Value2 = "https://www.upwork.com/jobs/mudded-video-editing_%7E01f58b6b0d38bf9a23?source=rss"
# Column I: text "$100" (currency-looking text -> force text format)
$ws.Cells.Item(133, 9).NumberFormat = "@"
$ws.Cells.Item(133, 9).Value2 = '$100'
$ws.Cells.Item(133, 10).Value2 = "June 15, 2024 09:04 UTC"
$ws.Cells.Item(133, 11).Value2 = "Video Editing"
$ws.Cells.Item(133, 12).Value2 = "Adobe After Effects,     Video Editing,     Video Post-Editing"
$ws.Cells.Item(133, 13).Value2 = "Spain"

# --- New row 134 data ---
$D134 = @'
I want to make one or a few Explainer Video(s) this weekend. I already have video footage of me (made in a recording studio) explaining the features and benefits of the product/service, so your work will be to edit this and combine it with visuals, graphics etc to make it into a compelling short video. 
Posted On
: June 15, 2024 08:41 UTC
Category
: Video Editing
Skills
:Video Editing,     Video Production    
Skills
:        Video Editing,                     Video Production            
Country
: United Kingdom
click to apply

'@
$E134 = @'
I want to make one or a few Explainer Video(s) this weekend. I already have video footage of me (made in a recording studio) explaining the features and benefits of the product/service, so your work will be to edit this and combine it with visuals, graphics etc to make it into a compelling short video. <br /><br /><br /><b>Posted On</b>: June 15, 2024 08:41 UTC<br /><b>Category</b>: Video Editing<br /><b>Skills</b>:Video Editing,     Video Production    
<br /><b>Skills</b>:        Video Editing,                     Video Production            <br /><b>Country</b>: United Kingdom
<br /><a href="https://www.upwork.com/jobs/Make-Explainer-Video-this-weekend_%7E019c213f850046d0c3?source=rss">click to apply</a>

'@

$ws.Cells.Item(134, 1).NumberFormat = "@"
$ws.Cells.Item(134, 1).Value2 = "23"
$ws.Cells.Item(134, 2).Value2 = "Make Explainer Video this weekend - Upwork"
$ws.Cells.Item(134, 3).Value2 = "https://www.upwork.com/jobs/Make-Explainer-Video-this-weekend_%7E019c213f850046d0c3?source=rss"
$ws.Cells.Item(134, 4).Value2 = $D134
$ws.Cells.Item(134, 5).Value2 = $E134
$ws.Cells.Item(134, 6).Value2 = "Sat, 15 Jun 2024 08:41:14 +0000"
$ws.Cells.Item(134, 7).Value2 = "https://www.upwork.com/jobs/Make-Explainer-Video-this-weekend_%7E019c213f850046d0c3?source=rss"
$ws.Cells.Item(134, 10).Value2 = "June 15, 2024 08:41 UTC"
$ws.Cells.Item(134, 11).Value2 = "Video Editing"
$ws.Cells.Item(134, 12).Value2 = "Video Editing,     Video Production"
$ws.Cells.Item(134, 13).Value2 = "United Kingdom"

Write-Host "Done"
